{"js": "// Rename the REST endpoint paths from camelCase to snake_case, e.g.\n// \"url: <SERVER>:80/addUser\" -> \"url: <SERVER>:80/add_user\"\n// (matches the document's actual textual change; the diff's\n// <w:proofErr> / run-splitting churn is cosmetic Word spell-check\n// bookkeeping, not a content edit).\n\nconst renames = [\n  [\"url: <SERVER>:80/addUser\", \"url: <SERVER>:80/add_user\"],\n  [\"url: <SERVER>:80/getItems\", \"url: <SERVER>:80/get_items\"],\n  [\"url: <SERVER>:80/likeItem\", \"url: <SERVER>:80/like_item\"],\n  [\"url: <SERVER>:80/getLikedItems\", \"url: <SERVER>:80/get_liked_items\"],\n  [\"url: <SERVER>:80/deleteLikedItem\", \"url: <SERVER>:80/delete_liked_item\"],\n];\n\nfor (const [oldText, newText] of renames) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Rename the REST endpoint paths from camelCase to snake_case, e.g.\n# \"url: <SERVER>:80/addUser\" -> \"url: <SERVER>:80/add_user\"\n# (matches the document's actual textual change; the diff's\n# <w:proofErr> / run-splitting churn is cosmetic Word spell-check\n# bookkeeping, not a content edit).\n\n$d = $word.ActiveDocument\n\n$renames = @{\n    \"url: <SERVER>:80/addUser\"         = \"url: <SERVER>:80/add_user\"\n    \"url: <SERVER>:80/getItems\"        = \"url: <SERVER>:80/get_items\"\n    \"url: <SERVER>:80/likeItem\"        = \"url: <SERVER>:80/like_item\"\n    \"url: <SERVER>:80/getLikedItems\"   = \"url: <SERVER>:80/get_liked_items\"\n    \"url: <SERVER>:80/deleteLikedItem\" = \"url: <SERVER>:80/delete_liked_item\"\n}\n\nforeach ($old in $renames.Keys) {\n    $new = $renames[$old]\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null  # wdReplaceAll\n}\n"}
